# dev brach change 1
# Append a new paragraph "IN DEV BRANCH" right after the existing
# "HELLO WORLD" paragraph, at the end of the document body.

$d = $word.ActiveDocument

# Collapse to the very end of the document story and insert a new,
# empty paragraph there (mirrors pressing Enter at the end of the doc).
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

# The freshly inserted paragraph is now the last paragraph in the
# document; fill it in with the new text.
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "IN DEV BRANCH"
